$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("TaskList")
$ws.Rows("33:33").Insert()
$ws.Range("A33").Value = "pager"
$ws.Range("B33").Value = "//*[@class='android.support.v4.view.ViewPager']"
$ws.Range("F62").Value = $ws.Range("F62").Value()
$ws.Range("F62").Value = ""
